$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4475.727
$ws.Range("I132").Value = 4613.8945
$ws.Range("K132").Value = 13841.6835
$ws.Range("M132").Value = -11311.6835
$ws.Range("H137").Value = 6921
$ws.Range("I137").Value = 6499.6665
$ws.Range("K137").Value = 19498.9995
$ws.Range("M137").Value = -16948.9995
$ws.Range("H138").Value = 6900.0815
$ws.Range("I138").Value = 6211.5
$ws.Range("J138").Value = 7123.4053
$ws.Range("K138").Value = 18634.5
$ws.Range("L138").Value = 21370.2159
$ws.Range("M138").Value = -13494.5
$ws.Range("N138").Value = -31650.2159

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20498.72
$ws.Range("I32").Value = 13887.622
$ws.Range("J32").Value = 79998.60000000001
$ws.Range("K32").Value = 13887.622
$ws.Range("L32").Value = 79998.60000000001
$ws.Range("M32").Value = -13600.622
$ws.Range("N32").Value = -80572.60000000001
$ws.Range("H132").Value = 15296988
$ws.Range("I132").Value = 16453922
$ws.Range("K132").Value = 49361766
$ws.Range("M132").Value = -49359236

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14572.444
$ws.Range("I20").Value = 4384.6665
$ws.Range("J20").Value = 34948
$ws.Range("K20").Value = 4384.6665
$ws.Range("L20").Value = 34948
$ws.Range("M20").Value = -4137.6665
$ws.Range("N20").Value = -35442
$ws.Range("H64").Value = 7255.625
$ws.Range("J64").Value = 8192.857
$ws.Range("L64").Value = 8192.857
$ws.Range("N64").Value = -8642.857
$ws.Range("H67").Value = 7255.625
$ws.Range("J67").Value = 8192.857
$ws.Range("L67").Value = 8192.857
$ws.Range("N67").Value = -9752.857
$ws.Range("H134").Value = 4231.375
$ws.Range("I134").Value = 3881.111
$ws.Range("J134").Value = 5282.1665
$ws.Range("K134").Value = 11643.333
$ws.Range("L134").Value = 15846.4995
$ws.Range("M134").Value = -9108.332999999999
$ws.Range("N134").Value = -20916.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 41670068
$ws.Range("I86").Value = 50002910
$ws.Range("J86").Value = 5836.25
$ws.Range("K86").Value = 50002910
$ws.Range("L86").Value = 5836.25
$ws.Range("M86").Value = -50001787
$ws.Range("N86").Value = -8082.25
$ws.Range("H89").Value = 41670068
$ws.Range("I89").Value = 50002910
$ws.Range("J89").Value = 5836.25
$ws.Range("K89").Value = 250014550
$ws.Range("L89").Value = 29181.25
$ws.Range("M89").Value = -250008934
$ws.Range("N89").Value = -40413.25
$ws.Range("H132").Value = 717803.1
$ws.Range("I132").Value = 3350
$ws.Range("J132").Value = 1432256.2
$ws.Range("K132").Value = 10050
$ws.Range("L132").Value = 4296768.6
$ws.Range("M132").Value = -7520
$ws.Range("N132").Value = -4301828.6
$ws.Range("H134").Value = 3439.0222
$ws.Range("I134").Value = 2925.9512
$ws.Range("J134").Value = 8698
$ws.Range("K134").Value = 8777.8536
$ws.Range("L134").Value = 26094
$ws.Range("M134").Value = -6242.8536
$ws.Range("N134").Value = -31164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2851329.8
$ws.Range("J113").Value = 3665424.2
$ws.Range("L113").Value = 10996272.6
$ws.Range("N113").Value = -11000612.6
$ws.Range("H131").Value = 32238.857
$ws.Range("J131").Value = 43739.8
$ws.Range("L131").Value = 131219.4
$ws.Range("N131").Value = -141299.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 178.15384
$ws.Range("I13").Value = 74.875
$ws.Range("K13").Value = 74.875
$ws.Range("M13").Value = 64.125
$ws.Range("H39").Value = 207765.25
$ws.Range("J39").Value = 207765.25
$ws.Range("L39").Value = 207765.25
$ws.Range("N39").Value = -208829.25
$ws.Range("H47").Value = 15333.333
$ws.Range("J47").Value = 15333.333
$ws.Range("L47").Value = 15333.333
$ws.Range("N47").Value = -16469.333
$ws.Range("H52").Value = 26250
$ws.Range("I52").Value = 15000
$ws.Range("J52").Value = 30000
$ws.Range("K52").Value = 15000
$ws.Range("L52").Value = 30000
$ws.Range("M52").Value = -14741
$ws.Range("N52").Value = -30518
$ws.Range("H80").Value = 71431430
$ws.Range("I80").Value = 2286.7144
$ws.Range("J80").Value = 142860580
$ws.Range("K80").Value = 2286.7144
$ws.Range("L80").Value = 142860580
$ws.Range("M80").Value = -1288.7144
$ws.Range("N80").Value = -142862576
$ws.Range("H83").Value = 71431430
$ws.Range("I83").Value = 2286.7144
$ws.Range("J83").Value = 142860580
$ws.Range("K83").Value = 11433.572
$ws.Range("L83").Value = 714302900
$ws.Range("M83").Value = -6441.572
$ws.Range("N83").Value = -714312884
$ws.Range("H126").Value = 7061.8687
$ws.Range("I126").Value = 9329
$ws.Range("K126").Value = 27987
$ws.Range("M126").Value = -25517
$ws.Range("H132").Value = 327636.75
$ws.Range("I132").Value = 403729.75
$ws.Range("J132").Value = 10582.5
$ws.Range("K132").Value = 1211189.25
$ws.Range("L132").Value = 31747.5
$ws.Range("M132").Value = -1208659.25
$ws.Range("N132").Value = -36807.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1412.5
$ws.Range("J22").Value = 2350
$ws.Range("L22").Value = 2350
$ws.Range("N22").Value = -2940
$ws.Range("H27").Value = 1412.5
$ws.Range("J27").Value = 2350
$ws.Range("L27").Value = 2350
$ws.Range("N27").Value = -2564
$ws.Range("H40").Value = 5150.647
$ws.Range("I40").Value = 4711.5
$ws.Range("J40").Value = 5390.1816
$ws.Range("K40").Value = 4711.5
$ws.Range("L40").Value = 5390.1816
$ws.Range("M40").Value = -4575.5
$ws.Range("N40").Value = -5662.1816
$ws.Range("H48").Value = 24999.666
$ws.Range("I48").Value = 24999.5
$ws.Range("K48").Value = 24999.5
$ws.Range("M48").Value = -24338.5
$ws.Range("H120").Value = 160000
$ws.Range("J120").Value = 160000
$ws.Range("L120").Value = 160000
$ws.Range("N120").Value = -169676
$ws.Range("H132").Value = 174520.69
$ws.Range("I132").Value = 273991.7
$ws.Range("J132").Value = 7228.5454
$ws.Range("K132").Value = 821975.1000000001
$ws.Range("L132").Value = 21685.6362
$ws.Range("M132").Value = -819445.1000000001
$ws.Range("N132").Value = -26745.6362
$ws.Range("H136").Value = 24396204
$ws.Range("I136").Value = 32263772
$ws.Range("K136").Value = 96791316
$ws.Range("M136").Value = -96788766

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 14999.5
$ws.Range("J38").Value = 14999.5
$ws.Range("L38").Value = 14999.5
$ws.Range("N38").Value = -15945.5
$ws.Range("H62").Value = 5481440.5
$ws.Range("I62").Value = 10954381
$ws.Range("K62").Value = 10954381
$ws.Range("M62").Value = -10953757
$ws.Range("H65").Value = 5481440.5
$ws.Range("I65").Value = 10954381
$ws.Range("K65").Value = 54771905
$ws.Range("M65").Value = -54768785
$ws.Range("H75").Value = 39666.332
$ws.Range("J75").Value = 39666.332
$ws.Range("L75").Value = 39666.332
$ws.Range("N75").Value = -41538.332
$ws.Range("H78").Value = 39666.332
$ws.Range("J78").Value = 39666.332
$ws.Range("L78").Value = 118998.996
$ws.Range("N78").Value = -128358.996
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H107").Value = 360.1111
$ws.Range("J107").Value = 305
$ws.Range("L107").Value = 915
$ws.Range("N107").Value = -4755
$ws.Range("H110").Value = 72000
$ws.Range("J110").Value = 72000
$ws.Range("L110").Value = 72000
$ws.Range("N110").Value = -80180
$ws.Range("H113").Value = 3371.7693
$ws.Range("I113").Value = 1185.5
$ws.Range("J113").Value = 6869.8
$ws.Range("K113").Value = 3556.5
$ws.Range("L113").Value = 20609.4
$ws.Range("M113").Value = -1386.5
$ws.Range("N113").Value = -24949.4
$ws.Range("H126").Value = 5460.2
$ws.Range("I126").Value = 3000.6667
$ws.Range("J126").Value = 9149.5
$ws.Range("K126").Value = 9002.000100000001
$ws.Range("L126").Value = 27448.5
$ws.Range("M126").Value = -6532.000100000001
$ws.Range("N126").Value = -32388.5
$ws.Range("H132").Value = 290410.38
$ws.Range("I132").Value = 307601.9
$ws.Range("K132").Value = 922805.7000000001
$ws.Range("M132").Value = -920275.7000000001
$ws.Range("H136").Value = 9638.607
$ws.Range("I136").Value = 10140.958
$ws.Range("J136").Value = 6624.5
$ws.Range("K136").Value = 30422.874
$ws.Range("L136").Value = 19873.5
$ws.Range("M136").Value = -27872.874
$ws.Range("N136").Value = -24973.5
